# Insert a new weekly price record for Albahaca (Hortaliza) at row 76.
# Excel shifts the existing rows 76-114 down to 77-115 (dimension grows
# from A1:R114 to A1:R115) and the new row 76 receives the fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 76, pushing rows 76..114 -> 77..115.
$ws.Rows.Item(76).Insert()

# Fill the newly inserted row 76 with the new observation.
$ws.Cells.Item(76, 1).Value = 8
$ws.Cells.Item(76, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44806
$ws.Cells.Item(76, 5).Value = 4
$ws.Cells.Item(76, 6).Value = 100112052
$ws.Cells.Item(76, 7).Value = "Albahaca"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 1000
$ws.Cells.Item(76, 11).Value = 4000
$ws.Cells.Item(76, 12).Value = 4500
$ws.Cells.Item(76, 13).Value = 4250
$ws.Cells.Item(76, 14).Value = "$/paquete"
$ws.Cells.Item(76, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(76, 16).Value = 4250
$ws.Cells.Item(76, 17).Value = 1
$ws.Cells.Item(76, 18).Value = "Hortaliza"
